# PSP_Sheet_2조.xlsx - "Add files via upload" update
# Updates the Time Recording Log (Sheet "작성자명") with additional
# logged sessions / corrected durations and replaces the two placeholder
# rows (13 & 14) with real dated entries, plus a brand-new row 15.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Row 9 (2019-09-24): stop time pushed back an hour, delta time +60 ---
$ws.Range("C9").Value = 0.91666666666666663
$ws.Range("E9").Value = 180

# --- Row 10 (2019-09-25): stop time +2h, interruption +40, delta +80 ---
$ws.Range("C10").Value = 0.95833333333333337
$ws.Range("D10").Value = 40
$ws.Range("E10").Value = 200

# --- Row 13: was a placeholder "10월 1일" text row -> real dated entry ---
$ws.Range("A13").Value = 43739
$ws.Range("B13").Value = 0.79166666666666663
$ws.Range("C13").Value = 0.875
$ws.Range("D13").Value = 0
$ws.Range("E13").Value = 120
$ws.Range("F13").Value = "Use Case Outline 관련 교수님 피드백 사항 점검 및 수정"

# --- Row 14: was a placeholder "10월 5일" text row -> real dated entry ---
$ws.Range("A14").Value = 43743
$ws.Range("B14").Value = 0.66666666666666663
$ws.Range("C14").Value = 0.77083333333333337
$ws.Range("D14").Value = 0
$ws.Range("E14").Value = 150
$ws.Range("F14").Value = "User 선정, 핵심로직 관련 구체화(시간표 구성 요건, 피드백로직)"

# --- Row 15: previously empty -> new logged session ---
$ws.Range("A15").Value = 43745
$ws.Range("B15").Value = 0.79166666666666663
$ws.Range("C15").Value = 0.91666666666666663
$ws.Range("D15").Value = 0
$ws.Range("E15").Value = 180
$ws.Range("F15").Value = "Use Case Spec 작성"

# --- View state: scrolled down one row, B15 now the active cell ---
$ws.Application.ActiveWindow.ScrollRow = 2
$ws.Range("B15").Select()
